$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1876923076923077
$ws.Range("C2").Value = 0.5661538461538461
$ws.Range("J2").Value = 0.009230769230769232
$ws.Range("P2").Value = 0.1292307692307692
$ws.Range("S2").Value = 0.1076923076923077
$ws.Range("B3").Value = 0.01041666666666667
$ws.Range("C3").Value = 0.02083333333333333
$ws.Range("J3").Value = 0.046875
$ws.Range("P3").Value = 0.765625
$ws.Range("S3").Value = 0.15625
$ws.Range("B6").Value = 0.05213270142180094
$ws.Range("D6").Value = 0.01895734597156398
$ws.Range("E6").Value = 0.004739336492890996
$ws.Range("F6").Value = 0.07109004739336493
$ws.Range("J6").Value = 0.2417061611374408
$ws.Range("O6").Value = 0.04265402843601896
$ws.Range("Q6").Value = 0.1658767772511848
$ws.Range("R6").Value = 0.06635071090047394
$ws.Range("S6").Value = 0.3364928909952606
$ws.Range("B7").Value = 0.1170731707317073
$ws.Range("D7").Value = 0.00975609756097561
$ws.Range("F7").Value = 0.02439024390243903
$ws.Range("J7").Value = 0.1170731707317073
$ws.Range("O7").Value = 0.02439024390243903
$ws.Range("Q7").Value = 0.1414634146341463
$ws.Range("R7").Value = 0.0975609756097561
$ws.Range("S7").Value = 0.4682926829268293
$ws.Range("B8").Value = 0.081374321880651
$ws.Range("D8").Value = 0.01989150090415913
$ws.Range("F8").Value = 0.05605786618444846
$ws.Range("J8").Value = 0.108499095840868
$ws.Range("O8").Value = 0.01265822784810127
$ws.Range("Q8").Value = 0.162748643761302
$ws.Range("R8").Value = 0.09584086799276673
$ws.Range("S8").Value = 0.4629294755877034
$ws.Range("B9").Value = 0.1437908496732026
$ws.Range("D9").Value = 0.0261437908496732
$ws.Range("F9").Value = 0.0457516339869281
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("O9").Value = 0.0261437908496732
$ws.Range("Q9").Value = 0.1699346405228758
$ws.Range("R9").Value = 0.1241830065359477
$ws.Range("S9").Value = 0.3464052287581699
$ws.Range("B10").Value = 0.1235431235431235
$ws.Range("D10").Value = 0.01476301476301476
$ws.Range("E10").Value = 0.000777000777000777
$ws.Range("F10").Value = 0.05905205905205906
$ws.Range("J10").Value = 0.1064491064491064
$ws.Range("O10").Value = 0.01476301476301476
$ws.Range("Q10").Value = 0.216006216006216
$ws.Range("R10").Value = 0.08547008547008547
$ws.Range("S10").Value = 0.3791763791763792
$ws.Range("G11").Value = 0.1784702549575071
$ws.Range("J11").Value = 0.1189801699716714
$ws.Range("K11").Value = 0.2379603399433428
$ws.Range("L11").Value = 0.4560906515580737
$ws.Range("S11").Value = 0.0084985835694051
$ws.Range("G12").Value = 0.7272727272727273
$ws.Range("J12").Value = 0.2181818181818182
$ws.Range("K12").Value = 0.01818181818181818
$ws.Range("L12").Value = 0.02424242424242424
$ws.Range("S12").Value = 0.01212121212121212
$ws.Range("G13").Value = 0.6086956521739131
$ws.Range("J13").Value = 0.3478260869565217
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.03349282296650718
$ws.Range("H15").Value = 0.2009569377990431
$ws.Range("I15").Value = 0.05741626794258373
$ws.Range("J15").Value = 0.3014354066985646
$ws.Range("K15").Value = 0.09090909090909091
$ws.Range("O15").Value = 0.0430622009569378
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("F16").Value = 0.02358490566037736
$ws.Range("H16").Value = 0.2075471698113208
$ws.Range("I16").Value = 0.07075471698113207
$ws.Range("J16").Value = 0.330188679245283
$ws.Range("K16").Value = 0.1320754716981132
$ws.Range("M16").Value = 0.02830188679245283
$ws.Range("O16").Value = 0.07075471698113207
$ws.Range("S16").Value = 0.1367924528301887
$ws.Range("F17").Value = 0.01545253863134658
$ws.Range("H17").Value = 0.1832229580573951
$ws.Range("I17").Value = 0.06622516556291391
$ws.Range("J17").Value = 0.4194260485651214
$ws.Range("K17").Value = 0.1169977924944812
$ws.Range("M17").Value = 0.01545253863134658
$ws.Range("N17").Value = 0.002207505518763797
$ws.Range("O17").Value = 0.04856512141280353
$ws.Range("S17").Value = 0.1324503311258278
$ws.Range("F18").Value = 0.02790697674418605
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.07906976744186046
$ws.Range("J18").Value = 0.386046511627907
$ws.Range("K18").Value = 0.08372093023255814
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("O18").Value = 0.06511627906976744
$ws.Range("S18").Value = 0.1348837209302326
$ws.Range("F19").Value = 0.01243599122165326
$ws.Range("H19").Value = 0.2523774689100219
$ws.Range("I19").Value = 0.0592538405267008
$ws.Range("J19").Value = 0.3679590343818581
$ws.Range("K19").Value = 0.1068032187271397
$ws.Range("M19").Value = 0.02194586686174104
$ws.Range("N19").Value = 0.000731528895391368
$ws.Range("O19").Value = 0.06071689831748354
$ws.Range("S19").Value = 0.1177761521580102
